$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlContinuous=1 (thin line style), xlLineStyleNone=-4142, xlCenter=-4108

# --- Update existing cell B5 text in place ---
$ws.Range("B5").Value = "return hello3() ;"

# --- Row 8: merged banner row (same visual style as rows 3 and 5) ---
$r8 = $ws.Range("B8:D8")
$r8.Merge()
$r8.Borders.LineStyle = 1
$r8.HorizontalAlignment = -4108
$ws.Range("B8").Value = "Rules String hello3()"

# --- Row 9: RET1 | (blank) | (blank), bordered table row ---
$b9 = $ws.Range("B9")
$b9.Borders.LineStyle = 1
$b9.Borders.Item(10).LineStyle = -4142
$b9.HorizontalAlignment = -4108

$c9 = $ws.Range("C9")
$c9.Borders.LineStyle = 1
$c9.Borders.Item(10).LineStyle = -4142
$c9.Borders.Item(7).LineStyle = -4142
$c9.HorizontalAlignment = -4108

$d9 = $ws.Range("D9")
$d9.Borders.LineStyle = 1
$d9.Borders.Item(7).LineStyle = -4142
$d9.HorizontalAlignment = -4108

$ws.Range("B9").Value = "RET1"
$ws.Range("B9:D9").Merge()

# --- Row 10: res ---
$b10 = $ws.Range("B10")
$b10.Borders.LineStyle = 1
$b10.Borders.Item(10).LineStyle = -4142
$b10.HorizontalAlignment = -4108

$c10 = $ws.Range("C10")
$c10.Borders.LineStyle = 1
$c10.Borders.Item(10).LineStyle = -4142
$c10.Borders.Item(7).LineStyle = -4142
$c10.HorizontalAlignment = -4108

$d10 = $ws.Range("D10")
$d10.Borders.LineStyle = 1
$d10.Borders.Item(7).LineStyle = -4142
$d10.HorizontalAlignment = -4108

$ws.Range("B10").Value = "res"
$ws.Range("B10:D10").Merge()

# --- Row 11: String res ---
$b11 = $ws.Range("B11")
$b11.Borders.LineStyle = 1
$b11.Borders.Item(10).LineStyle = -4142
$b11.HorizontalAlignment = -4108

$c11 = $ws.Range("C11")
$c11.Borders.LineStyle = 1
$c11.Borders.Item(10).LineStyle = -4142
$c11.Borders.Item(7).LineStyle = -4142
$c11.HorizontalAlignment = -4108

$d11 = $ws.Range("D11")
$d11.Borders.LineStyle = 1
$d11.Borders.Item(7).LineStyle = -4142
$d11.HorizontalAlignment = -4108

$ws.Range("B11").Value = "String res"
$ws.Range("B11:D11").Merge()

# --- Row 12: From dependency ---
$b12 = $ws.Range("B12")
$b12.Borders.LineStyle = 1
$b12.Borders.Item(10).LineStyle = -4142
$b12.HorizontalAlignment = -4108

$c12 = $ws.Range("C12")
$c12.Borders.LineStyle = 1
$c12.Borders.Item(10).LineStyle = -4142
$c12.Borders.Item(7).LineStyle = -4142
$c12.HorizontalAlignment = -4108

$d12 = $ws.Range("D12")
$d12.Borders.LineStyle = 1
$d12.Borders.Item(7).LineStyle = -4142
$d12.HorizontalAlignment = -4108

$ws.Range("B12").Value = "From dependency"
$ws.Range("B12:D12").Merge()

# --- Row 13: formula-looking text, needs quote prefix ---
$b13 = $ws.Range("B13")
$b13.Borders.LineStyle = 1
$b13.Borders.Item(10).LineStyle = -4142
$b13.HorizontalAlignment = -4108

$c13 = $ws.Range("C13")
$c13.Borders.LineStyle = 1
$c13.Borders.Item(10).LineStyle = -4142
$c13.Borders.Item(7).LineStyle = -4142
$c13.HorizontalAlignment = -4108

$d13 = $ws.Range("D13")
$d13.Borders.LineStyle = 1
$d13.Borders.Item(7).LineStyle = -4142
$d13.HorizontalAlignment = -4108

$ws.Range("B13").Value = "'=return helloFromDependency() ;"
$ws.Range("B13:D13").Merge()

# --- Row 15: Environment label ---
$r15 = $ws.Range("B15:C15")
$r15.HorizontalAlignment = -4108
$ws.Range("B15").Value = "Environment"
$r15.Merge()

# --- Row 16: dependency / Module3_2 ---
$ws.Range("B16").Value = "dependency"
$ws.Range("C16").Value = "Module3_2"

# --- View state ---
$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
